$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.881.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.120.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.117.69'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("E9").Value = '  +3.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("E11").Value = '  -1.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.409'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.652.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.137'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.837.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.139.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '350.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("E26").Value = '  -1.46%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0872'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.31%  '
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("E34").Value = '  -5.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  -2.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.69'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0669'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.01'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.695'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.397.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '36.93'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.155.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0266'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.960'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.739'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.57%  '
